$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = 112083118
$ws.Range("B7").Value = 94134
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 53
$ws.Range("F7").Value = "Vedtrappmossa"
$ws.Range("G7").Value = "Crossocalyx hellerianus"
$ws.Range("H7").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q7").Value = 412576.6879626553
$ws.Range("R7").Value = 6656303.56951345

# Row 8
$ws.Range("A8").Value = 112083111
$ws.Range("B8").Value = 90666
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 4364
$ws.Range("F8").Value = "Dropptaggsvamp"
$ws.Range("G8").Value = "Hydnellum ferrugineum"
$ws.Range("H8").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q8").Value = 412204.6634863199
$ws.Range("R8").Value = 6655988.977203708

# Row 9
$ws.Range("A9").Value = 112083125
$ws.Range("B9").Value = 89369
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 5447
$ws.Range("F9").Value = "Vedticka"
$ws.Range("G9").Value = "Fuscoporia viticola"
$ws.Range("H9").Value = "(Schwein.) Murrill"
$ws.Range("Q9").Value = 413015.9403039298
$ws.Range("R9").Value = 6656414.640994807

# Row 10
$ws.Range("A10").Value = 112083128
$ws.Range("B10").Value = 77186
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 353
$ws.Range("F10").Value = "Dvärgbägarlav"
$ws.Range("G10").Value = "Cladonia parasitica"
$ws.Range("H10").Value = "(Hoffm.) Hoffm."
$ws.Range("Q10").Value = 413190.1061828797
$ws.Range("R10").Value = 6656475.01450387

# Row 11
$ws.Range("A11").Value = 112083126
$ws.Range("B11").Value = 78536
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 229497
$ws.Range("F11").Value = "Korallblylav"
$ws.Range("G11").Value = "Parmeliella triptophylla"
$ws.Range("H11").Value = "(Ach.) Müll.Arg."
$ws.Range("Q11").Value = 413016.7201701452
$ws.Range("R11").Value = 6656341.641577623

# Row 12
$ws.Range("A12").Value = 112083127
$ws.Range("B12").Value = 77604
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 6450
$ws.Range("F12").Value = "Skuggblåslav"
$ws.Range("G12").Value = "Hypogymnia vittata"
$ws.Range("H12").Value = "(Ach.) Parrique"
$ws.Range("Q12").Value = 413051.8096683071
$ws.Range("R12").Value = 6656343.312587639

# Row 13
$ws.Range("A13").Value = 112083110
$ws.Range("B13").Value = 78107
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 6453
$ws.Range("F13").Value = "Vedskivlav"
$ws.Range("G13").Value = "Hertelidea botryosa"
$ws.Range("H13").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q13").Value = 412205.6393663768
$ws.Range("R13").Value = 6656050.944565876

# Row 14
$ws.Range("A14").Value = 112083112
$ws.Range("B14").Value = 79444
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 1049
$ws.Range("F14").Value = "Kortskaftad ärgspik"
$ws.Range("G14").Value = "Microcalicium ahlneri"
$ws.Range("H14").Value = "Tibell"
$ws.Range("Q14").Value = 412283.7604491137
$ws.Range("R14").Value = 6656072.080045181
